$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.021.68'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '1.861.06'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('D5').Value = '312.19'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('D7').Value = '0.5098'
$ws.Range('E7').Value = '  +1.47%  '
$ws.Range('D8').Value = '0.3845'
$ws.Range('E8').Value = '  +0.62%  '
$ws.Range('D9').Value = '0.08294'
$ws.Range('E9').Value = '  -7.23%  '
$ws.Range('D10').Value = '1.113'
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('D12').Value = '6.229'
$ws.Range('E12').Value = '  -2.14%  '
$ws.Range('D13').Value = '20.58'
$ws.Range('D14').Value = '1.861.90'
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('D15').Value = '7.227'
$ws.Range('D16').Value = '1.004'
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('D17').Value = '0.00001097'
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('D18').Value = '90.86'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('D19').Value = '0.06636'
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').Value = '17.69'
$ws.Range('E20').Value = '  -2.65%  '
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').Value = "'6.040"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.26%  '
$ws.Range('D23').Value = '28.049.36'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('E24').Value = '  -3.29%  '
$ws.Range('D25').Value = '2.228'
$ws.Range('E25').Value = '  -1.06%  '
$ws.Range('D26').Value = '2.539'
$ws.Range('E26').Value = '  +1.89%  '
$ws.Range('D27').Value = '2.074.03'
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('D28').Value = '157.96'
$ws.Range('E28').Value = '  +0.40%  '
$ws.Range('D29').Value = '20.54'
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('D30').Value = '124.91'
$ws.Range('E30').Value = '  -1.13%  '
$ws.Range('E31').Value = '  -0.86%  '
$ws.Range('D32').Value = '1.038'
$ws.Range('E32').Value = '  -1.40%  '
$ws.Range('D33').Value = '5.841'
$ws.Range('E33').Value = '  +4.43%  '
$ws.Range('D34').Value = '3.601'
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('D35').Value = '9.431'
$ws.Range('E35').Value = '  +0.16%  '
$ws.Range('D36').Value = "'0.02420"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.95%  '
$ws.Range('D37').Value = "'0.06530"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.84%  '
$ws.Range('D38').Value = '0.2172'
$ws.Range('E38').Value = '  -0.61%  '
$ws.Range('D39').Value = '1.204'
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('D40').Value = '0.6466'
$ws.Range('E40').Value = '  +1.54%  '
$ws.Range('D41').Value = '1.225'
$ws.Range('E41').Value = '  -4.56%  '
$ws.Range('D42').Value = '4.978'
$ws.Range('E42').Value = '  +1.74%  '
$ws.Range('D43').Value = '11.19'
$ws.Range('E43').Value = '  -2.22%  '
$ws.Range('D44').Value = '0.6103'
$ws.Range('E44').Value = '  +1.56%  '
$ws.Range('D45').Value = '13.07'
$ws.Range('E45').Value = '  -0.80%  '
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('D47').Value = '3.659'
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('D48').Value = '2.015'
$ws.Range('E48').Value = '  +1.12%  '
$ws.Range('D49').Value = '1.209'
$ws.Range('E49').Value = '  -1.98%  '
$ws.Range('D50').Value = '120.02'
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('D51').Value = '78.33'
$ws.Range('E51').Value = '  -1.02%  '
